$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serials (Excel 1900 date system) for Dec-31 of each year 2004-2024,
# replacing the previous "YYYYQ4" text labels in column A (rows 2-22).
$dates = @(38352, 38717, 39082, 39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735, 43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
